$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.09179766666666667
$ws.Range("H2").Value = 0.275393
$ws.Range("M2").Value = 5.828378333333333
$ws.Range("N2").Value = 17.485135
$ws.Range("O2").Value = 0.1413867973615592
$ws.Range("P2").Value = 0.1413867973615592
$ws.Range("Q2").Value = 0.5350315314505556
$ws.Range("R2").Value = 4.815283783054999
$ws.Range("S2").Value = 0.1413867973615592
$ws.Range("T2").Value = 0.1413867973615592

# Row 3
$ws.Range("G3").Value = 0.09179766666666667
$ws.Range("H3").Value = 0.275393
$ws.Range("O3").Value = 0.604557320991465
$ws.Range("P3").Value = 0.604557320991465
$ws.Range("Q3").Value = 2.287747055140889
$ws.Range("R3").Value = 20.589723496268
$ws.Range("S3").Value = 0.604557320991465
$ws.Range("T3").Value = 0.604557320991465

# Row 4
$ws.Range("G4").Value = 0.09179766666666667
$ws.Range("H4").Value = 0.275393
$ws.Range("M4").Value = 10.47292833333333
$ws.Range("N4").Value = 31.418785
$ws.Range("O4").Value = 0.2540558816469758
$ws.Range("P4").Value = 0.2540558816469758
$ws.Range("Q4").Value = 0.9613903841672222
$ws.Range("R4").Value = 8.652513457505
$ws.Range("S4").Value = 0.2540558816469758
$ws.Range("T4").Value = 0.2540558816469758
